$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: change from "silver" to "black"
$ws.Range("B3").Value = "black"

# B4: change quantity from 4 to 1
$ws.Range("B4").Value = 1

# B10: change from a plain value to a formula computing weight * 1.053
$ws.Range("B10").Formula = '=B4*1.053'

# B13: change formula to VLOOKUP-based price lookup
$ws.Range("B13").Formula = '=VLOOKUP($B$1,$D$2:$E$5,2,FALSE)'

# Remove row 14 entirely (it contained the old VLOOKUP formula, now moved to B13)
$ws.Range("A14:H14").Delete()
